$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 11 through 26 currently have a row height of 19.5; set them to 18.75
$ws.Range("A11:A26").EntireRow.RowHeight = 18.75
